$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "Oula Fuel Marketing Company - KSCP (KWSE:OULAFUEL)"

# Row 2
$ws.Range("D2").Value = 0.0345
$ws.Range("E2").Value = -0.11
$ws.Range("G2").Value = 0.005485854858548586
$ws.Range("H2").Value = 0.005485854858548586
$ws.Range("I2").Value = -0.02600246002460025
$ws.Range("J2").Value = -0.02575628319539214
$ws.Range("K2").Value = 13.3
$ws.Range("L2").Value = 0.01635916359163592
$ws.Range("M2").Value = 0.202
$ws.Range("N2").Value = 0.0006198220312979443
$ws.Range("O2").Value = 0.01518796992481203
$ws.Range("P2").Value = 0.202
$ws.Range("Q2").Value = 0.0006198220312979443
$ws.Range("R2").Value = 0.01518796992481203
$ws.Range("U2").Value = 68
$ws.Range("V2").Value = 0.2086529610309911
$ws.Range("W2").Value = 0.03054104100308408
$ws.Range("X2").Value = 0.06738611127674589
$ws.Range("Y2").Value = -0.03684507027366181
$ws.Range("Z2").Value = 1.745009658725048
$ws.Range("AA2").Value = -0.04791554433705814
$ws.Range("AB2").Value = 0.05478006596602802
$ws.Range("AC2").Value = -0.1026956103030862
$ws.Range("AD2").Value = 120
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 120
$ws.Range("AG2").Value = 52
$ws.Range("AH2").Value = 0.2691186364655753
$ws.Range("AI2").Value = 0.2228412256267409
$ws.Range("AJ2").Value = 0.1376025403545912
$ws.Range("AK2").Value = 0.1105207226354942
$ws.Range("AL2").Value = 1.29
$ws.Range("AM2").Value = -4.33
$ws.Range("AN2").Value = -85.1063829787234
$ws.Range("AO2").Value = -16.38759689922481
$ws.Range("AP2").Value = -36.87943262411348
$ws.Range("AQ2").Value = 4.882217090069284

# Row 3
$ws.Range("D3").Value = 0.0347
$ws.Range("E3").Value = -0.116
$ws.Range("G3").Value = 0.002674846625766871
$ws.Range("H3").Value = 0.002674846625766871
$ws.Range("I3").Value = -0.02292024539877301
$ws.Range("J3").Value = -0.02273403956466185
$ws.Range("K3").Value = 6.71
$ws.Range("L3").Value = 0.01646625766871166
$ws.Range("M3").Value = 0.202
$ws.Range("N3").Value = 0.001234718826405868
$ws.Range("O3").Value = 0.0301043219076006
$ws.Range("P3").Value = 0.202
$ws.Range("Q3").Value = 0.001234718826405868
$ws.Range("R3").Value = 0.0301043219076006
$ws.Range("U3").Value = 45.4
$ws.Range("V3").Value = 0.2775061124694376
$ws.Range("W3").Value = 0.02942982456140351
$ws.Range("X3").Value = 0.0811903595298402
$ws.Range("Y3").Value = -0.0517605349684367
$ws.Range("Z3").Value = 1.456397426733381
$ws.Range("AA3").Value = -0.0331097967212284
$ws.Range("AB3").Value = 0.05597826890840447
$ws.Range("AC3").Value = -0.08908806562963287
$ws.Range("AD3").Value = 120
$ws.Range("AF3").Value = 120
$ws.Range("AG3").Value = 74.59999999999999
$ws.Range("AH3").Value = 0.4231311706629055
$ws.Range("AI3").Value = 0.3642987249544627
$ws.Range("AJ3").Value = 0.3131821998320739
$ws.Range("AK3").Value = 0.2626760563380282
$ws.Range("AL3").Value = 1.29
$ws.Range("AM3").Value = -0.54
$ws.Range("AN3").Value = 12000
$ws.Range("AO3").Value = -7.240310077519379
$ws.Range("AP3").Value = 7459.999999999999
$ws.Range("AQ3").Value = 17.29629629629629

# Row 4
$ws.Range("D4").Value = 0.0343
$ws.Range("E4").Value = -0.104
$ws.Range("G4").Value = 0.008310727496917386
$ws.Range("H4").Value = 0.008310727496917386
$ws.Range("I4").Value = -0.02909987669543773
$ws.Range("J4").Value = -0.02878528343386543
$ws.Range("K4").Value = 6.59
$ws.Range("L4").Value = 0.01625154130702836
$ws.Range("M4").Value = -0
$ws.Range("N4").Value = -0
$ws.Range("O4").Value = -0
$ws.Range("P4").Value = -0
$ws.Range("Q4").Value = -0
$ws.Range("R4").Value = -0
$ws.Range("U4").Value = 22.6
$ws.Range("V4").Value = 0.1392483056069008
$ws.Range("W4").Value = 0.03165225744476465
$ws.Range("X4").Value = 0.05358186302365157
$ws.Range("Y4").Value = -0.02192960557888692
$ws.Range("Z4").Value = 2.178936055883933
$ws.Range("AA4").Value = -0.06272129195288788
$ws.Range("AB4").Value = 0.05358186302365157
$ws.Range("AC4").Value = -0.1163031549765395
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 0
$ws.Range("AG4").Value = -22.6
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = -0.1617752326413744
$ws.Range("AK4").Value = -0.1211796246648794
$ws.Range("AM4").Value = -3.79
$ws.Range("AN4").Value = -0
$ws.Range("AP4").Value = 15.91549295774648
$ws.Range("AQ4").Value = 3.113456464379948

$ws.Range("T4").ClearContents()